$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20, which shifts the existing rows 20-25
# (weekly price records) down to rows 21-26.
$ws.Rows("20:20").Insert()

# Populate the new row 20 with the new weekly record.
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44466
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 100112026
$ws.Range("G20").Value = "Haba"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 9000
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Región de O'Higgins"
$ws.Range("P20").Value = 360
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
